# FEATS_brm.xls work: add handle/metal/ceramic/plastic/glass/steel concept
# columns, fix a couple of stray values, and append the new concept rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- small corrections to existing data -----------------------------------
$ws.Range("S6").Value  = "gray"            # fork color was "silver" -> "gray"
$ws.Range("I10").Value = "steel"           # spatula material "stainless_steel" -> "steel"

# --- new "Concepts" rows for the new feature columns ----------------------
$ws.Range("A12").Value = "handle"
$ws.Range("A13").Value = "metal"
$ws.Range("A14").Value = "ceramic"
$ws.Range("A15").Value = "plastic"
$ws.Range("A16").Value = "glass"
$ws.Range("A17").Value = "steel"

# --- handle_feature / handle_val (columns T:U) ----------------------------
$ws.Range("T2").Value = "handle_feature"
$ws.Range("U2").Value = "handle_val"
$ws.Range("T3").Value = "shape"
$ws.Range("U3").Value = "long"
$ws.Range("T4").Value = "material"
$ws.Range("U4").Value = "plastic"
$ws.Range("T5").Value = "material"
$ws.Range("U5").Value = "metal"
$ws.Range("T6").Value = "material"
$ws.Range("U6").Value = "wood"

# --- metal_feature / metal_val (columns V:W) ------------------------------
$ws.Range("V2").Value = "metal_feature"
$ws.Range("W2").Value = "metal_val"
$ws.Range("V3").Value = "is"
$ws.Range("W3").Value = "material"
$ws.Range("V4").Value = "color"
$ws.Range("W4").Value = "gray"
$ws.Range("V5").Value = "color"
$ws.Range("W5").Value = "different_colors"

# --- ceramic_feature / ceramic_val (columns X:Y) --------------------------
$ws.Range("X2").Value = "ceramic_feature"
$ws.Range("Y2").Value = "ceramic_val"
$ws.Range("X3").Value = "is"
$ws.Range("Y3").Value = "material"
$ws.Range("X4").Value = "color"
$ws.Range("Y4").Value = "white"
$ws.Range("X5").Value = "color"
$ws.Range("Y5").Value = "different_colors"

# --- plastic_feature / plastic_val (columns Z:AA) -------------------------
$ws.Range("Z2").Value  = "plastic_feature"
$ws.Range("AA2").Value = "plastic_val"
$ws.Range("Z3").Value  = "is"
$ws.Range("AA3").Value = "material"
$ws.Range("Z4").Value  = "color"
$ws.Range("AA4").Value = "different_colors"

# --- glass_feature / glass_val (columns AB:AC) ----------------------------
$ws.Range("AB2").Value = "glass_feature"
$ws.Range("AC2").Value = "glass_val"
$ws.Range("AB3").Value = "is"
$ws.Range("AC3").Value = "material"
$ws.Range("AB4").Value = "texture"
$ws.Range("AC4").Value = "shiny"
$ws.Range("AB5").Value = "color"
$ws.Range("AC5").Value = "transparent"

# --- steel_feature / steel_val (columns AD:AE) ----------------------------
$ws.Range("AD2").Value = "steel_feature"
$ws.Range("AE2").Value = "steel_val"
$ws.Range("AD3").Value = "is"
$ws.Range("AE3").Value = "metal"
$ws.Range("AD4").Value = "color"
$ws.Range("AE4").Value = "gray"
$ws.Range("AD5").Value = "texture"
$ws.Range("AE5").Value = "shiny"

# --- column widths for the new columns (autofit-style, matches the
#     widths Excel settled on for columns T,U,W,X,Y,Z,AA) -----------------
$ws.Columns.Item(20).ColumnWidth = 12.333333333333332  # T  -> 13.109375
$ws.Columns.Item(21).ColumnWidth = 10.666666666666668  # U  -> 11.5546875
$ws.Columns.Item(23).ColumnWidth = 13.166666666666668  # W  -> 14
$ws.Columns.Item(24).ColumnWidth = 13.333333333333332  # X  -> 14.21875
$ws.Columns.Item(25).ColumnWidth = 12.0                # Y  -> 12.77734375
$ws.Columns.Item(26).ColumnWidth = 13.5                # Z  -> 14.33203125
$ws.Columns.Item(27).ColumnWidth = 11.833333333333332  # AA -> 12.6640625

# --- leave the selection where the author left it -------------------------
$ws.Range("A2").Select()
